$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-Rows {
    param($Worksheet, $Rows, $Columns)

    $n = $Rows.Count

    # Snapshot current values for each row/column before any writes happen
    $snapshot = @{}
    foreach ($r in $Rows) {
        $rowVals = @{}
        foreach ($c in $Columns) {
            $addr = "$c$r"
            $rowVals[$c] = $Worksheet.Range($addr).Value2
        }
        $snapshot[$r] = $rowVals
    }

    # Each row takes on the values that were previously held by the row
    # above it in the group, with the first row wrapping around to take
    # the last row's original values.
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $Rows[$i]
        $srcIndex = ($i - 1 + $n) % $n
        $srcRow = $Rows[$srcIndex]
        foreach ($c in $Columns) {
            $addr = "$c$destRow"
            $Worksheet.Range($addr).Value2 = $snapshot[$srcRow][$c]
        }
    }
}

$columns = @("B", "D", "E", "F", "G")

Rotate-Rows $ws @(149, 150) $columns
Rotate-Rows $ws @(161, 162, 163) $columns
Rotate-Rows $ws @(264, 265) $columns
Rotate-Rows $ws @(279, 280) $columns
Rotate-Rows $ws @(351, 352) $columns
Rotate-Rows $ws @(372, 373) $columns
Rotate-Rows $ws @(379, 380) $columns
Rotate-Rows $ws @(457, 458) $columns
Rotate-Rows $ws @(536, 537) $columns
Rotate-Rows $ws @(583, 584) $columns
Rotate-Rows $ws @(586, 587) $columns
Rotate-Rows $ws @(590, 591) $columns
Rotate-Rows $ws @(593, 594) $columns
Rotate-Rows $ws @(601, 602) $columns
Rotate-Rows $ws @(687, 688) $columns
Rotate-Rows $ws @(889, 890) $columns
